# Tasks.xlsx update: "Data splits and Pickling"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark a few WEEK 1 tasks as Completed in the Status column (D)
$ws.Range("D9").Value2 = "Completed"
$ws.Range("D10").Value2 = "Completed"
$ws.Range("D13").Value2 = "Completed"

# Owner of "Pickling of Data" (row 10) changes from Sayantika to Abhijit
$ws.Range("C10").Value2 = "Abhijit"

# Grow the table by one row and add the new WEEK 2 task
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Range("A21").Value2 = "WEEK 2"
$ws.Range("B21").Value2 = "Improvement of visualizations"
$ws.Range("C21").Value2 = "Sayantika"
$ws.Range("A21").HorizontalAlignment = -4108

# Move the active selection to D13, matching the saved workbook state
$ws.Range("D13").Select() | Out-Null
